$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 103.4275383333333
$ws.Range("H2").Value = 310.282615
$ws.Range("I2").Value = 0.2485530285127421
$ws.Range("J2").Value = 0.2485530285127421
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.628421
$ws.Range("N2").Value = 4.885263
$ws.Range("O2").Value = 0.048329411442081
$ws.Range("P2").Value = 0.048329411442081
$ws.Range("Q2").Value = 168.423575400305
$ws.Range("R2").Value = 1515.812178602745
$ws.Range("S2").Value = 0.0120124215801676
$ws.Range("T2").Value = 0.0120124215801676

# Row 3
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 103.4275383333333
$ws.Range("H3").Value = 310.282615
$ws.Range("I3").Value = 0.2485530285127421
$ws.Range("J3").Value = 0.2485530285127421
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.459065000000001
$ws.Range("N3").Value = 28.377195
$ws.Range("O3").Value = 0.2807327123897247
$ws.Range("P3").Value = 0.2807327123897247
$ws.Range("Q3").Value = 978.3278078849919
$ws.Range("R3").Value = 8804.950270964926
$ws.Range("S3").Value = 0.06977696586706268
$ws.Range("T3").Value = 0.06977696586706267

# Row 4
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 103.4275383333333
$ws.Range("H4").Value = 310.282615
$ws.Range("I4").Value = 0.2485530285127421
$ws.Range("J4").Value = 0.2485530285127421
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.6418243333333334
$ws.Range("N4").Value = 1.925473
$ws.Range("O4").Value = 0.01904850912583786
$ws.Range("P4").Value = 0.01904850912583786
$ws.Range("Q4").Value = 66.38231083909946
$ws.Range("R4").Value = 597.4407975518951
$ws.Range("S4").Value = 0.004734564631879606
$ws.Range("T4").Value = 0.004734564631879605

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 103.4275383333333
$ws.Range("H5").Value = 310.282615
$ws.Range("I5").Value = 0.2485530285127421
$ws.Range("J5").Value = 0.2485530285127421
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.96489266666667
$ws.Range("N5").Value = 65.894678
$ws.Range("O5").Value = 0.6518893670423564
$ws.Range("P5").Value = 0.6518893670423563
$ws.Range("Q5").Value = 2271.774778269219
$ws.Range("R5").Value = 20445.97300442297
$ws.Range("S5").Value = 0.1620290764336322
$ws.Range("T5").Value = 0.1620290764336322

# Row 6
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 216.130539
$ws.Range("H6").Value = 648.391617
$ws.Range("I6").Value = 0.5193964865470273
$ws.Range("J6").Value = 0.5193964865470272
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 1.628421
$ws.Range("N6").Value = 4.885263
$ws.Range("O6").Value = 0.048329411442081
$ws.Range("P6").Value = 0.048329411442081
$ws.Range("Q6").Value = 351.951508448919
$ws.Range("R6").Value = 3167.563576040271
$ws.Range("S6").Value = 0.02510212649990257
$ws.Range("T6").Value = 0.02510212649990257

# Row 7
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 216.130539
$ws.Range("H7").Value = 648.391617
$ws.Range("I7").Value = 0.5193964865470273
$ws.Range("J7").Value = 0.5193964865470272
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 9.459065000000001
$ws.Range("N7").Value = 28.377195
$ws.Range("O7").Value = 0.2807327123897247
$ws.Range("P7").Value = 0.2807327123897247
$ws.Range("Q7").Value = 2044.392816886035
$ws.Range("R7").Value = 18399.53535197431
$ws.Range("S7").Value = 0.1458115844740402
$ws.Range("T7").Value = 0.1458115844740401

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 216.130539
$ws.Range("H8").Value = 648.391617
$ws.Range("I8").Value = 0.5193964865470273
$ws.Range("J8").Value = 0.5193964865470272
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.6418243333333334
$ws.Range("N8").Value = 1.925473
$ws.Range("O8").Value = 0.01904850912583786
$ws.Range("P8").Value = 0.01904850912583786
$ws.Range("Q8").Value = 138.717839106649
$ws.Range("R8").Value = 1248.460551959841
$ws.Range("S8").Value = 0.009893728713919171
$ws.Range("T8").Value = 0.009893728713919168

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 216.130539
$ws.Range("H9").Value = 648.391617
$ws.Range("I9").Value = 0.5193964865470273
$ws.Range("J9").Value = 0.5193964865470272
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 21.96489266666667
$ws.Range("N9").Value = 65.894678
$ws.Range("O9").Value = 0.6518893670423564
$ws.Range("P9").Value = 0.6518893670423563
$ws.Range("Q9").Value = 4747.284091123815
$ws.Range("R9").Value = 42725.55682011433
$ws.Range("S9").Value = 0.3385890468591654
$ws.Range("T9").Value = 0.3385890468591653

# Row 10
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 71.607325
$ws.Range("H10").Value = 214.821975
$ws.Range("I10").Value = 0.1720839321833696
$ws.Range("J10").Value = 0.1720839321833696
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 1.628421
$ws.Range("N10").Value = 4.885263
$ws.Range("O10").Value = 0.048329411442081
$ws.Range("P10").Value = 0.048329411442081
$ws.Range("Q10").Value = 116.606871783825
$ws.Range("R10").Value = 1049.461846054425
$ws.Range("S10").Value = 0.008316715161061233
$ws.Range("T10").Value = 0.008316715161061233

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 71.607325
$ws.Range("H11").Value = 214.821975
$ws.Range("I11").Value = 0.1720839321833696
$ws.Range("J11").Value = 0.1720839321833696
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 9.459065000000001
$ws.Range("N11").Value = 28.377195
$ws.Range("O11").Value = 0.2807327123897247
$ws.Range("P11").Value = 0.2807327123897247
$ws.Range("Q11").Value = 677.3383416511251
$ws.Range("R11").Value = 6096.045074860125
$ws.Range("S11").Value = 0.04830958904052678
$ws.Range("T11").Value = 0.04830958904052678

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 71.607325
$ws.Range("H12").Value = 214.821975
$ws.Range("I12").Value = 0.1720839321833696
$ws.Range("J12").Value = 0.1720839321833696
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.6418243333333334
$ws.Range("N12").Value = 1.925473
$ws.Range("O12").Value = 0.01904850912583786
$ws.Range("P12").Value = 0.01904850912583786
$ws.Range("Q12").Value = 45.95932362990834
$ws.Range("R12").Value = 413.6339126691751
$ws.Range("S12").Value = 0.003277942352604978
$ws.Range("T12").Value = 0.003277942352604978

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 71.607325
$ws.Range("H13").Value = 214.821975
$ws.Range("I13").Value = 0.1720839321833696
$ws.Range("J13").Value = 0.1720839321833696
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 21.96489266666667
$ws.Range("N13").Value = 65.894678
$ws.Range("O13").Value = 0.6518893670423564
$ws.Range("P13").Value = 0.6518893670423563
$ws.Range("Q13").Value = 1572.847207772117
$ws.Range("R13").Value = 14155.62486994905
$ws.Range("S13").Value = 0.1121796856291766
$ws.Range("T13").Value = 0.1121796856291766

# Row 14 (new)
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "App"
$ws.Range("C14").Value = "Gpc1"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 24.953198
$ws.Range("H14").Value = 74.859594
$ws.Range("I14").Value = 0.05996655275686102
$ws.Range("J14").Value = 0.05996655275686102
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 1.628421
$ws.Range("N14").Value = 4.885263
$ws.Range("O14").Value = 0.048329411442081
$ws.Range("P14").Value = 0.048329411442081
$ws.Range("Q14").Value = 40.63431164035801
$ws.Range("R14").Value = 365.708804763222
$ws.Range("S14").Value = 0.002898148200949593
$ws.Range("T14").Value = 0.002898148200949593

# Row 15 (new)
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "App"
$ws.Range("C15").Value = "Gpc1"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 24.953198
$ws.Range("H15").Value = 74.859594
$ws.Range("I15").Value = 0.05996655275686102
$ws.Range("J15").Value = 0.05996655275686102
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 9.459065000000001
$ws.Range("N15").Value = 28.377195
$ws.Range("O15").Value = 0.2807327123897247
$ws.Range("P15").Value = 0.2807327123897247
$ws.Range("Q15").Value = 236.03392183987
$ws.Range("R15").Value = 2124.30529655883
$ws.Range("S15").Value = 0.01683457300809512
$ws.Range("T15").Value = 0.01683457300809512

# Row 16 (new)
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "App"
$ws.Range("C16").Value = "Gpc1"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 24.953198
$ws.Range("H16").Value = 74.859594
$ws.Range("I16").Value = 0.05996655275686102
$ws.Range("J16").Value = 0.05996655275686102
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.6418243333333334
$ws.Range("N16").Value = 1.925473
$ws.Range("O16").Value = 0.01904850912583786
$ws.Range("P16").Value = 0.01904850912583786
$ws.Range("Q16").Value = 16.01556967088467
$ws.Range("R16").Value = 144.140127037962
$ws.Range("S16").Value = 0.001142273427434105
$ws.Range("T16").Value = 0.001142273427434104

# Row 17 (new)
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "App"
$ws.Range("C17").Value = "Gpc1"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 24.953198
$ws.Range("H17").Value = 74.859594
$ws.Range("I17").Value = 0.05996655275686102
$ws.Range("J17").Value = 0.05996655275686102
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 21.96489266666667
$ws.Range("N17").Value = 65.894678
$ws.Range("O17").Value = 0.6518893670423564
$ws.Range("P17").Value = 0.6518893670423563
$ws.Range("Q17").Value = 548.0943157600814
$ws.Range("R17").Value = 4932.848841840732
$ws.Range("S17").Value = 0.03909155812038221
$ws.Range("T17").Value = 0.0390915581203822
